$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.043.50"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "2.888.04"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("E4").Value = "  +0.04%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "351.13"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.36%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "111.37"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.61%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.558"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("E8").Value = "  +0.05%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.622"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.16%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.81"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("E11").Value = "  +0.38%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0857"
$c.ClearFormats()
$ws.Range("E12").Value = "  +2.08%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "19.92"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.57%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.76"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").Value = "3.344.37"
$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.915.98"
$ws.Range("E16").Value = "  +4.22%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E17").Value = "  +5.80%  "

$ws.Range("D18").Value = "52.090.08"
$ws.Range("E18").Value = "  +0.85%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.97%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.ClearFormats()
$ws.Range("E20").Value = "  +4.31%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.49"
$c.ClearFormats()
$ws.Range("E21").Value = "  +8.01%  "

$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +0.63%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "70.68"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.22%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "269.34"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.36%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.78"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.49"
$c.ClearFormats()
$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("E28").Value = "  -0.18%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.51"
$c.ClearFormats()
$ws.Range("E29").Value = "  +1.35%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.36"
$c.ClearFormats()
$ws.Range("E30").Value = "  +2.30%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.24"
$c.ClearFormats()
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.43"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.49%  "

$ws.Range("E33").Value = "  +8.44%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0943"
$c.ClearFormats()
$ws.Range("E34").Value = "  +10.06%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "52.85"
$c.ClearFormats()
$ws.Range("E35").Value = "  +1.32%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0457"
$c.ClearFormats()
$ws.Range("E36").Value = "  +2.19%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.04%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.ClearFormats()
$ws.Range("E38").Value = "  +4.62%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.57"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.92%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.ClearFormats()
$ws.Range("E40").Value = "  +2.26%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.64"
$c.ClearFormats()
$ws.Range("E41").Value = "  +5.96%  "

$ws.Range("E42").Value = "  +1.59%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "22.71"
$c.ClearFormats()
$ws.Range("E43").Value = "  +3.65%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "122.09"
$c.ClearFormats()
$ws.Range("E44").Value = "  +1.68%  "

$ws.Range("E45").Value = "  +0.53%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.ClearFormats()
$ws.Range("E46").Value = "  +3.64%  "

$ws.Range("D47").Value = "2.194.56"
$ws.Range("E47").Value = "  +2.48%  "

$ws.Range("E48").Value = "  +5.87%  "

$ws.Range("E49").Value = "  +23.21%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.942"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.70%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0324"
$c.ClearFormats()
$ws.Range("E51").Value = "  +10.28%  "
